$d = $word.ActiveDocument
$full = $d.Content.Text
$idx = $full.IndexOf("Exercise:  ")
$end = $idx + ("Exercise:  ").Length
$insPos = $idx + 8

$rSuffix = $d.Range($insPos, $end)
$rSuffix.Text = ""
$p1 = $d.Range($insPos, $insPos)
$p1.InsertAfter(" 1")

$tmpRange = $d.Range($insPos, $insPos)
$d.Bookmarks.Add("zzTempSplit", $tmpRange)
$d.Bookmarks("zzTempSplit").Delete()

Write-Host "Checkpoint B:" $d.Content.Text.Substring($idx, 11)

# position right after " 1" = insPos+2
$afterSpace1 = $insPos + 2
$p2 = $d.Range($afterSpace1, $afterSpace1)
$p2.InsertBefore(":  ")

Write-Host "Checkpoint C:" $d.Content.Text.Substring($idx, 14)
